# Replace the comparative-analysis table contents with the new
# "select metrics for selected companies" page.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (overwrite text in place so existing header styling/fill/border survives)
$ws.Range("B1").Value = "Debt equity ratio"
$ws.Range("C1").Value = "Interest coverage"
$ws.Range("D1").Value = "Revenue growth"
$ws.Range("E1").Value = "Retained earnings"

# Row 2 - AAPL
$ws.Range("A2").Value = "AAPL"
$ws.Range("B2").Value = 1.787532584558942
$ws.Range("C2").Value = 29.06203915586067
$ws.Range("D2").Value = -0.02041077580526742
$ws.Range("E2").Value = -214000000

# Row 3 - SONY
$ws.Range("A3").Value = "SONY"
$ws.Range("B3").Value = 0.5620591368200297
$ws.Range("C3").Value = 20.49508914182966
$ws.Range("D3").Value = 0.014244529073212
$ws.Range("E3").Value = 4614637

# Old table had 5 data rows (2-5) and 6 columns (A-F); new table only needs
# rows 1-3 and columns A-E, so drop the now-unused rows/column.
$ws.Range("A4:F5").EntireRow.Delete() | Out-Null
$ws.Range("F1:F3").EntireColumn.Delete() | Out-Null

# The percent-style formatting that used to highlight the "Price to book
# ratio" column (old col E) now highlights "Interest coverage" (col C);
# reset the old column back to the default, unformatted style.
$ws.Range("E2:E3").ClearFormats()
$ws.Range("C2:C3").NumberFormat = "0.00%"
